$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# so numeric-looking strings (e.g. "509.68") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '58.579.84'
$ws.Range("E2").Value = '  -4.19%  '
$ws.Range("D3").Value = '2.557.48'
$ws.Range("E3").Value = '  -3.87%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '509.68'
$ws.Range("E5").Value = '  -4.55%  '
$ws.Range("D6").Value = '146.49'
$ws.Range("E6").Value = '  -6.51%  '
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -4.03%  '
$ws.Range("D9").Value = '2.571.10'
$ws.Range("E9").Value = '  -3.84%  '
$ws.Range("D10").Value = '6.24'
$ws.Range("E10").Value = '  -5.19%  '
$ws.Range("E11").Value = '  -6.23%  '
$ws.Range("E12").Value = '  -5.03%  '
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '3.004.40'
$ws.Range("E14").Value = '  -3.76%  '
$ws.Range("D15").Value = '58.546.66'
$ws.Range("E15").Value = '  -4.19%  '
$ws.Range("D16").Value = '21.03'
$ws.Range("E16").Value = '  -4.86%  '
$ws.Range("D17").Value = '0.0000136'
$ws.Range("E17").Value = '  -5.30%  '
$ws.Range("D18").Value = '2.567.05'
$ws.Range("E18").Value = '  -3.70%  '
$ws.Range("D19").Value = '347.37'
$ws.Range("E19").Value = '  -2.19%  '
$ws.Range("D21").Value = '10.24'
$ws.Range("E21").Value = '  -4.42%  '
$ws.Range("E22").Value = '  -4.37%  '
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("D24").Value = '60.64'
$ws.Range("E24").Value = '  -1.64%  '
$ws.Range("D25").Value = '0.414'
$ws.Range("E25").Value = '  -4.48%  '
$ws.Range("D26").Value = '0.997'
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("E27").Value = '  -5.23%  '
$ws.Range("D28").Value = '2.666.95'
$ws.Range("E28").Value = '  -3.72%  '
$ws.Range("D29").Value = '0.0₃0801'
$ws.Range("E29").Value = '  -7.22%  '
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -5.06%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  -0.06%  '
$ws.Range("E32").Value = '  -3.65%  '
$ws.Range("E33").Value = '  -4.99%  '
$ws.Range("D34").Value = '149.61'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("E35").Value = '  -6.11%  '
$ws.Range("E36").Value = '  -4.78%  '
$ws.Range("D37").Value = '0.900'
$ws.Range("E37").Value = '  +1.42%  '
$ws.Range("D38").Value = '1.13'
$ws.Range("E38").Value = '  -6.36%  '
$ws.Range("E39").Value = '  -8.91%  '
$ws.Range("D40").Value = '35.99'
$ws.Range("E40").Value = '  -2.33%  '
$ws.Range("E41").Value = '  -6.40%  '
$ws.Range("D42").Value = '287.08'
$ws.Range("E42").Value = '  -6.73%  '
$ws.Range("D43").Value = '3.57'
$ws.Range("E43").Value = '  -6.76%  '
$ws.Range("E44").Value = '  -2.75%  '
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("E46").Value = '  -6.97%  '
$ws.Range("D47").Value = '0.0537'
$ws.Range("E47").Value = '  -5.04%  '
$ws.Range("D48").Value = '19.08'
$ws.Range("E48").Value = '  -5.48%  '
$ws.Range("E49").Value = '  -4.39%  '
$ws.Range("D50").Value = '10.26'
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").Value = '4.62'
$ws.Range("E51").Value = '  -7.72%  '
